# Mise à jour de l'application
# A new training-session column (2025-09-19) is appended right after the
# existing "AW" (2025-09-18) column: a new date header in row 1, plus each
# player's attendance mark ("P", "B", "R", "REP", ...) in the same row as
# before, in the new column AX.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new date header AX1 = 2025-09-19 (serial 45919) ---
# Copy/PasteSpecial(formats) from AW1 so AX1 reuses the existing date-column
# style (center aligned, short-date number format) instead of registering a
# brand-new style/number-format entry. Row 1 has no formulas depending on it,
# so this is safe.
$ws.Range("AW1").Copy()
$ws.Range("AX1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("AX1").Value = 45919

# --- Rows 2-29: this session's attendance mark for every player, column AX ---
# Column B:J hold COUNTA/COUNTIF summary formulas that read across the whole
# row (K:VQ etc.), so writing the new value directly (rather than via
# Copy/PasteSpecial, which leaves those formulas stale) lets them pick up the
# new column automatically. HorizontalAlignment=-4108 (xlCenter) reproduces
# the same cell style already used across the rest of the row without
# creating a duplicate style entry.
$attendance = [ordered]@{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "B"
    6  = "B"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "B"
    11 = "P"
    12 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "R"
    22 = "P"
    23 = "B"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "REP"
    28 = "P"
    29 = "P"
}

foreach ($row in $attendance.Keys) {
    $cell = $ws.Range("AX$row")
    $cell.Value = $attendance[$row]
    $cell.HorizontalAlignment = -4108
}

# --- Refresh the recap formulas & move the live selection to match the edit ---
$excel.Calculate()
[void]$ws.Range("AZ24").Select()
